$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.023.21"
$ws.Range("E2").Value = "  -1.23%  "

$ws.Range("D3").Value = "3.847.57"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.17"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("D7").Value = "3.843.66"
$ws.Range("E7").Value = "  -2.01%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.49%  "

$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").Value = "4.486.29"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "3.828.90"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("D17").Value = "68.026.29"
$ws.Range("E17").Value = "  -1.46%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("E24").Value = "  -3.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.99%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("E30").Value = "  -0.91%  "

$ws.Range("D31").Value = "3.991.05"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.75"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.74%  "

$ws.Range("E33").Value = "  -3.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.26%  "

$ws.Range("D35").Value = "3.815.24"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +9.50%  "

$ws.Range("E40").Value = "  -3.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("E42").Value = "  -2.98%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "425.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.14%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000271"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.85%  "

$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
